$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C8').Value = 56
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.0'
$ws.Range('E8').Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('C9').Value = 32
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '4'
$ws.Range('E9').Value = 'Long point  (up to 10 mtr.)'
$ws.Range('F9').Value = 662
$ws.Range('G9').NumberFormat = "@"
$ws.Range('G9').Value = '21184.00'
$ws.Range('A10').Value = ''
$ws.Range('C10').Value = 26
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '2.0'
$ws.Range('E10').Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F10').Value = 0
$ws.Range('G10').NumberFormat = "@"
$ws.Range('G10').Value = '0.00'
$ws.Range('C11').Value = 41
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '4.0'
$ws.Range('E11').Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F11').Value = 50
$ws.Range('G11').NumberFormat = "@"
$ws.Range('G11').Value = '2050.00'
$ws.Range('C12').Value = 51
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.0'
$ws.Range('E12').Value = 'Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F12').Value = 33
$ws.Range('G12').NumberFormat = "@"
$ws.Range('G12').Value = '1683.00'
$ws.Range('C13').Value = 11
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.0'
$ws.Range('E13').Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F13').Value = 303
$ws.Range('G13').NumberFormat = "@"
$ws.Range('G13').Value = '3333.00'
$ws.Range('A14').Value = 'R. mtr.'
$ws.Range('C14').Value = 84
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '17'
$ws.Range('E14').Value = '25 mm'
$ws.Range('F14').Value = 56
$ws.Range('G14').NumberFormat = "@"
$ws.Range('G14').Value = '4704.00'
$ws.Range('A15').Value = 'Set'
$ws.Range('C15').Value = 58
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '13.0'
$ws.Range('E15').Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range('F15').Value = 5733
$ws.Range('G15').NumberFormat = "@"
$ws.Range('G15').Value = '332514.00'
$ws.Range('C16').Value = 14
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '23'
$ws.Range('E16').Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range('F16').Value = 20
$ws.Range('G16').NumberFormat = "@"
$ws.Range('G16').Value = '280.00'
$ws.Range('A17').Value = ''
$ws.Range('C17').Value = 2
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '15.0'
$ws.Range('E17').Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range('F17').Value = 0
$ws.Range('G17').NumberFormat = "@"
$ws.Range('G17').Value = '0.00'
$ws.Range('C18').Value = 87
$ws.Range('G18').NumberFormat = "@"
$ws.Range('G18').Value = '164430.00'
$ws.Range('C19').Value = 61
$ws.Range('A20').Value = ''
$ws.Range('C20').Value = 71
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.0'
$ws.Range('E20').Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range('F20').Value = 0
$ws.Range('G20').NumberFormat = "@"
$ws.Range('G20').Value = '0.00'
$ws.Range('A21').Value = 'Each'
$ws.Range('C21').Value = 36
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '30'
$ws.Range('E21').Value = ' 6 A to 32 A rating'
$ws.Range('F21').Value = 187
$ws.Range('G21').NumberFormat = "@"
$ws.Range('G21').Value = '6732.00'
$ws.Range('C22').Value = 75
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '31'
$ws.Range('E22').Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range('C23').Value = 85
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '35'
$ws.Range('E23').Value = '8 Way (8+2)'
$ws.Range('F23').Value = 2184
$ws.Range('G23').NumberFormat = "@"
$ws.Range('G23').Value = '185640.00'
$ws.Range('A24').Value = '%'
$ws.Range('C24').Value = 51
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '37'
$ws.Range('E24').Value = 'Add Tender Premium '

# Remove the old "Grand Total" item row (row 25); this shifts the summary
# rows (Grand Total Rs. / Tender Premium / NET PAYABLE AMOUNT) up by one.
$ws.Rows.Item(25).Delete()

# Update the summary amounts to reflect the new Grand Total.
$ws.Range('G26').NumberFormat = "@"
$ws.Range('G26').Value = '722550.00'
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H26').Value = '722550.00'

$ws.Range('G28').NumberFormat = "@"
$ws.Range('G28').Value = '722550.00'
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H28').Value = '722550.00'
